$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.991.67"
$ws.Range("E2").Value = "  -1.76%  "

# Row 3
$ws.Range("D3").Value = "3.243.21"
$ws.Range("E3").Value = "  -0.98%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.44%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.51%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "3.239.67"
$ws.Range("E9").Value = "  -1.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.47%  "

# Row 11
$ws.Range("E11").Value = "  +0.78%  "

# Row 12
$ws.Range("E12").Value = "  -3.07%  "

# Row 13
$ws.Range("D13").Value = "3.803.96"
$ws.Range("E13").Value = "  -1.19%  "

# Row 14
$ws.Range("E14").Value = "  -3.11%  "

# Row 15
$ws.Range("D15").Value = "65.046.60"
$ws.Range("E15").Value = "  -1.75%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.89"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.74%  "

# Rows 17 & 18: swap ShibaInu and WrappedEther, with updated price/volume
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000160"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.33%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.232.86"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "419.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.68%  "

# Row 20
$ws.Range("E20").Value = "  -2.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.90"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.23"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.34%  "

# Row 23
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.10"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.53%  "

# Row 26
$ws.Range("E26").Value = "  +4.53%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.497"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.81%  "

# Row 28
$ws.Range("E28").Value = "  -0.77%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.14%  "

# Row 30
$ws.Range("E30").Value = "  -0.19%  "

# Row 31
$ws.Range("E31").Value = "  -3.76%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.91"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.68%  "

# Row 33
$ws.Range("E33").Value = "  +0.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.03"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.08%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.46"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.08%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.41"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.79%  "

# Row 38
$ws.Range("E38").Value = "  -1.60%  "

# Row 39
$ws.Range("D39").Value = "2.838.56"
$ws.Range("E39").Value = "  +2.35%  "

# Row 40
$ws.Range("E40").Value = "  -2.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.58"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.26"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.55%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.726"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0632"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "304.20"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.22"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.25%  "

# Row 50
$ws.Range("E50").Value = "  -0.68%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.102"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.06%  "

